$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.330.28"
$ws.Range("E2").Value = "'  +1.95%  "

$ws.Range("D3").Value = "'2.392.07"
$ws.Range("E3").Value = "'  +7.54%  "

$ws.Range("E4").Value = "'  -0.26%  "

$ws.Range("D5").Value = "'323.59"
$ws.Range("E5").Value = "'  +11.59%  "

$ws.Range("D6").Value = "'105.73"
$ws.Range("E6").Value = "'  -5.01%  "

$ws.Range("D7").Value = "'0.655"
$ws.Range("E7").Value = "'  +4.74%  "

$ws.Range("E8").Value = "'  -0.01%  "

$ws.Range("D9").Value = "'0.655"
$ws.Range("E9").Value = "'  +9.27%  "

$ws.Range("D10").Value = "'41.98"
$ws.Range("E10").Value = "'  -4.21%  "

$ws.Range("D11").Value = "'0.0946"
$ws.Range("E11").Value = "'  +3.71%  "

$ws.Range("D12").Value = "'8.59"
$ws.Range("E12").Value = "'  -0.54%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'17.32"
$ws.Range("E13").Value = "'  +16.22%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'1.02"
$ws.Range("E14").Value = "'  +0.34%  "

$ws.Range("E15").Value = "'  +2.90%  "

$ws.Range("D16").Value = "'2.756.04"
$ws.Range("E16").Value = "'  +7.68%  "

$ws.Range("D17").Value = "'2.392.85"
$ws.Range("E17").Value = "'  +5.73%  "

$ws.Range("D18").Value = "'43.310.99"
$ws.Range("E18").Value = "'  +2.05%  "

$ws.Range("D19").Value = "'0.0000109"
$ws.Range("E19").Value = "'  +4.01%  "

$ws.Range("D20").Value = "'7.38"
$ws.Range("E20").Value = "'  +3.23%  "

$ws.Range("D21").Value = "'75.77"
$ws.Range("E21").Value = "'  +3.70%  "

$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").Value = "'3.46"
$ws.Range("E22").Value = "'  +3.93%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'268.14"
$ws.Range("E23").Value = "'  +13.72%  "

$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "'  +1.58%  "

$ws.Range("D25").Value = "'9.77"
$ws.Range("E25").Value = "'  +8.47%  "

$ws.Range("D26").Value = "'11.89"
$ws.Range("E26").Value = "'  +4.19%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  -0.12%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'22.86"
$ws.Range("E28").Value = "'  +7.85%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'177.52"
$ws.Range("E29").Value = "'  +2.43%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "'  +0.24%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.90"
$ws.Range("E31").Value = "'  +0.85%  "

$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "'  +3.77%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0930"
$ws.Range("E33").Value = "'  +5.76%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.90"
$ws.Range("E34").Value = "'  +5.07%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.134"
$ws.Range("E35").Value = "'  +6.74%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'4.86"
$ws.Range("E36").Value = "'  -2.49%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.15"
$ws.Range("E37").Value = "'  -0.83%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0369"
$ws.Range("E38").Value = "'  -2.57%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.110"
$ws.Range("E39").Value = "'  +4.84%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.83"
$ws.Range("E40").Value = "'  +18.43%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.58"
$ws.Range("E41").Value = "'  +20.72%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'126.10"
$ws.Range("E42").Value = "'  +24.94%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.233"
$ws.Range("E43").Value = "'  +1.17%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'69.47"
$ws.Range("E44").Value = "'  -3.24%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  +0.06%  "

$ws.Range("D46").Value = "'12.64"
$ws.Range("E46").Value = "'  +2.01%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.64"
$ws.Range("E47").Value = "'  +14.66%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'5.66"
$ws.Range("E48").Value = "'  +5.89%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'85.34"
$ws.Range("E49").Value = "'  +52.92%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.32"
$ws.Range("E50").Value = "'  +3.65%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.606.81"
$ws.Range("E51").Value = "'  +12.94%  "

# Remove the quotePrefix formatting artifact from the leading apostrophe
# used above to force text storage, restoring original (default) style.
$ws.Range("D2:E51").ClearFormats()
